$wb = $excel.ActiveWorkbook

# The two sheets "展览" and "全部类型" contain identical data tables and
# both need the same "想去人数" (F column) values updated.
$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2 = 2172
    3 = 1648
    5 = 1062
    6 = 629
    7 = 32
    8 = 5746
    9 = 84
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
